$d = $word.ActiveDocument

# The page used to end with a blank spacer paragraph followed by a
# "Ver no Jupiter Salvar em pdf Salvar em docx" line and a
# "© 2020 . Contact: ..." colophon line, right after the last bibliography
# paragraph ("...Basic Food Microbiology...USA,1970."). Those three
# paragraphs were dropped from the page footer while the trailing blank
# paragraph and the page-break paragraph after them were kept.
#
# Find the last bibliography paragraph dynamically (rather than hardcoding
# an index) so the edit is anchored to content, then remove the three
# paragraphs that immediately follow it.
$anchorIndex = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "*Basic Food Microbiology*") {
        $anchorIndex = $i
    }
}

if ($anchorIndex -ne $null) {
    $start = $d.Paragraphs($anchorIndex + 1).Range.Start
    $end = $d.Paragraphs($anchorIndex + 3).Range.End
    $d.Range($start, $end).Delete()
}
